$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sopot"
$ws.Range("C2").Value = 0.154
$ws.Range("B3").Value = "Kraków"
$ws.Range("C3").Value = 0.273
$ws.Range("B4").Value = "Warszawa"
$ws.Range("C4").Value = 0.294
$ws.Range("B5").Value = "Toruń"
$ws.Range("C5").Value = 0.296
$ws.Range("B6").Value = "Siedlce"
$ws.Range("C6").Value = 0.332
$ws.Range("B7").Value = "Poznań"
$ws.Range("C7").Value = 0.345
$ws.Range("B8").Value = "Rzeszów"
$ws.Range("C8").Value = 0.348
$ws.Range("B9").Value = "Nowy Sącz"
$ws.Range("C9").Value = 0.349
$ws.Range("B10").Value = "Jelenia Góra"
$ws.Range("C10").Value = 0.357
$ws.Range("B11").Value = "Grudziądz"
$ws.Range("C11").Value = 0.359
$ws.Range("B12").Value = "Świnoujście"
$ws.Range("C12").Value = 0.367
$ws.Range("B13").Value = "Olsztyn"
$ws.Range("C13").Value = 0.379
$ws.Range("B14").Value = "Zamość"
$ws.Range("C14").Value = 0.382
$ws.Range("B15").Value = "Koszalin"
$ws.Range("C15").Value = 0.383
$ws.Range("B16").Value = "Lublin"
$ws.Range("C16").Value = 0.384
$ws.Range("B17").Value = "Wrocław"
$ws.Range("C17").Value = 0.389
$ws.Range("B18").Value = "Elbląg"
$ws.Range("C18").Value = 0.403
$ws.Range("B19").Value = "Katowice"
$ws.Range("C19").Value = 0.404
$ws.Range("B20").Value = "Szczecin"
$ws.Range("C20").Value = 0.406
$ws.Range("B21").Value = "Skierniewice"
$ws.Range("C21").Value = 0.409
$ws.Range("B22").Value = "Jastrzębie-Zdrój"
$ws.Range("C22").Value = 0.414
$ws.Range("B23").Value = "Chorzów"
$ws.Range("C23").Value = 0.418
$ws.Range("B24").Value = "Gdańsk"
$ws.Range("C24").Value = 0.421
$ws.Range("B25").Value = "Jaworzno"
$ws.Range("C25").Value = 0.429
$ws.Range("B26").Value = "Zielona Góra"
$ws.Range("C26").Value = 0.43
$ws.Range("B27").Value = "Bielsko-Biała"
$ws.Range("C27").Value = 0.43
$ws.Range("B28").Value = "Kielce"
$ws.Range("C28").Value = 0.43
$ws.Range("B29").Value = "Płock"
$ws.Range("C29").Value = 0.434
$ws.Range("B30").Value = "Kalisz"
$ws.Range("C30").Value = 0.436
$ws.Range("B31").Value = "Ostrołęka"
$ws.Range("C31").Value = 0.437
$ws.Range("B32").Value = "Przemyśl"
$ws.Range("C32").Value = 0.437
$ws.Range("B33").Value = "Słupsk"
$ws.Range("C33").Value = 0.441
$ws.Range("B34").Value = "Bydgoszcz"
$ws.Range("C34").Value = 0.444
$ws.Range("B35").Value = "Legnica"
$ws.Range("C35").Value = 0.451
$ws.Range("B36").Value = "Krosno"
$ws.Range("C36").Value = 0.451
$ws.Range("B37").Value = "Radom"
$ws.Range("C37").Value = 0.454
$ws.Range("B38").Value = "Łódź"
$ws.Range("C38").Value = 0.457
$ws.Range("B39").Value = "Zabrze"
$ws.Range("C39").Value = 0.458
$ws.Range("B40").Value = "Suwałki"
$ws.Range("C40").Value = 0.461
$ws.Range("B41").Value = "Tarnów"
$ws.Range("C41").Value = 0.462
$ws.Range("B42").Value = "Leszno"
$ws.Range("C42").Value = 0.471
$ws.Range("B43").Value = "Wałbrzych"
$ws.Range("C43").Value = 0.474
$ws.Range("B44").Value = "Tarnobrzeg"
$ws.Range("C44").Value = 0.482
$ws.Range("B45").Value = "Opole"
$ws.Range("C45").Value = 0.487
$ws.Range("B46").Value = "Siemianowice Śląskie"
$ws.Range("C46").Value = 0.492
$ws.Range("B47").Value = "Łomża"
$ws.Range("C47").Value = 0.493
$ws.Range("B48").Value = "Gorzów Wielkopolski"
$ws.Range("C48").Value = 0.498
$ws.Range("B49").Value = "Konin"
$ws.Range("C49").Value = 0.501
$ws.Range("B50").Value = "Chełm"
$ws.Range("C50").Value = 0.504
$ws.Range("B51").Value = "Tychy"
$ws.Range("C51").Value = 0.51
$ws.Range("B52").Value = "Piotrków Trybunalski"
$ws.Range("C52").Value = 0.513
$ws.Range("B53").Value = "Gdynia"
$ws.Range("C53").Value = 0.513
$ws.Range("B54").Value = "Częstochowa"
$ws.Range("C54").Value = 0.513
$ws.Range("B55").Value = "Rybnik"
$ws.Range("C55").Value = 0.515
$ws.Range("B56").Value = "Bytom"
$ws.Range("C56").Value = 0.521
$ws.Range("B57").Value = "Gliwice"
$ws.Range("C57").Value = 0.521
$ws.Range("B58").Value = "Biała Podlaska"
$ws.Range("C58").Value = 0.526
$ws.Range("B59").Value = "Białystok"
$ws.Range("C59").Value = 0.531
$ws.Range("B60").Value = "Sosnowiec"
$ws.Range("C60").Value = 0.534
$ws.Range("B61").Value = "Świętochłowice"
$ws.Range("C61").Value = 0.538
$ws.Range("B62").Value = "Dąbrowa Górnicza"
$ws.Range("C62").Value = 0.547
$ws.Range("B63").Value = "Włocławek"
$ws.Range("C63").Value = 0.557
$ws.Range("B64").Value = "Mysłowice"
$ws.Range("C64").Value = 0.566
$ws.Range("B65").Value = "Piekary Śląskie"
$ws.Range("C65").Value = 0.605
$ws.Range("B66").Value = "Ruda Śląska"
$ws.Range("C66").Value = 0.607
$ws.Range("B67").Value = "Żory"
$ws.Range("C67").Value = 0.611

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C67"), 0, 1, 0, 0)
$ws.Sort.SetRange($ws.Range("B2:C67"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Range("J11").Select() | Out-Null
